$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update URL value
$ws.Range("B2").Value = "https://hl7.fr/fhir/fr/medication/ValueSet/FrMethodOfAdministration"

# Update Date value
$ws.Range("B8").Value = "2024-12-26T10:27:36+00:00"

# Clear the Copyright value cell (row 14, col B) -- text removed entirely
$ws.Range("B14").ClearContents()
